$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 361, shifting the existing rows 361..435
# down to 363..437 (dimension grows from R435 to R437).
$ws.Rows("361:362").Insert()

# Populate the newly inserted row 361 with its data.
$ws.Range("A361").Value = 5
$ws.Range("B361").Value = "Macroferia Regional de Talca"
$ws.Range("C361").Value = "Maule"
$ws.Range("D361").Value = 44641
$ws.Range("E361").Value = 7
$ws.Range("F361").Value = 100114001
$ws.Range("G361").Value = "Papa"
$ws.Range("H361").Value = "Asterix"
$ws.Range("I361").Value = "1a (cosecha)"
$ws.Range("J361").Value = 1600
$ws.Range("K361").Value = 7000
$ws.Range("L361").Value = 7000
$ws.Range("M361").Value = 7000
$ws.Range("N361").Value = '$/saco 25 kilos'
$ws.Range("O361").Value = "Región de Los Lagos"
$ws.Range("P361").Value = 280
$ws.Range("Q361").Value = 25
$ws.Range("R361").Value = "Hortaliza"

# Populate the newly inserted row 362 with its data.
$ws.Range("A362").Value = 5
$ws.Range("B362").Value = "Macroferia Regional de Talca"
$ws.Range("C362").Value = "Maule"
$ws.Range("D362").Value = 44641
$ws.Range("E362").Value = 7
$ws.Range("F362").Value = 100114001
$ws.Range("G362").Value = "Papa"
$ws.Range("H362").Value = "Rodeo"
$ws.Range("I362").Value = "1a (cosecha)"
$ws.Range("J362").Value = 1800
$ws.Range("K362").Value = 6500
$ws.Range("L362").Value = 6500
$ws.Range("M362").Value = 6500
$ws.Range("N362").Value = '$/saco 25 kilos'
$ws.Range("O362").Value = "Región de Los Lagos"
$ws.Range("P362").Value = 260
$ws.Range("Q362").Value = 25
$ws.Range("R362").Value = "Hortaliza"
